$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Format D2:E51 as text before writing so numeric-looking strings
# (e.g. "0.9977", "312.81") are stored as text, not auto-converted numbers,
# matching the source data which stores these as inline strings.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "27.543.92"
$ws.Range("E2").Value = "  -4.92%  "
$ws.Range("D3").Value = "1.841.25"
$ws.Range("E3").Value = "  -4.33%  "
$ws.Range("D4").Value = "0.9977"
$ws.Range("E4").Value = "  -0.74%  "
$ws.Range("D5").Value = "312.81"
$ws.Range("E5").Value = "  -3.78%  "
$ws.Range("D6").Value = "0.9989"
$ws.Range("E6").Value = "  -0.50%  "
$ws.Range("D7").Value = "0.4255"
$ws.Range("E7").Value = "  -7.25%  "
$ws.Range("D8").Value = "0.3622"
$ws.Range("E8").Value = "  -5.15%  "
$ws.Range("D9").Value = "43.43"
$ws.Range("E9").Value = "  -4.78%  "
$ws.Range("D10").Value = "0.07201"
$ws.Range("E10").Value = "  -7.08%  "
$ws.Range("D11").Value = "0.8978"
$ws.Range("E11").Value = "  -8.29%  "
$ws.Range("D12").Value = "20.67"
$ws.Range("E12").Value = "  -8.42%  "
$ws.Range("D13").Value = "1.847.65"
$ws.Range("E13").Value = "  -3.87%  "
$ws.Range("D14").Value = "6.588"
$ws.Range("E14").Value = "  -5.50%  "
$ws.Range("D15").Value = "5.320"
$ws.Range("E15").Value = "  -6.79%  "
$ws.Range("D16").Value = "0.06787"
$ws.Range("E16").Value = "  -2.83%  "
$ws.Range("D17").Value = "0.9993"
$ws.Range("E17").Value = "  -0.70%  "
$ws.Range("D18").Value = "77.25"
$ws.Range("E18").Value = "  -8.93%  "
$ws.Range("D19").Value = "0.000008944"
$ws.Range("E19").Value = "  -5.77%  "
$ws.Range("D20").Value = "0.9994"
$ws.Range("E20").Value = "  -0.40%  "
$ws.Range("D21").Value = "15.36"
$ws.Range("E21").Value = "  -7.97%  "
$ws.Range("D22").Value = "27.516.17"
$ws.Range("E22").Value = "  -5.08%  "
$ws.Range("D23").Value = "4.934"
$ws.Range("E23").Value = "  -7.69%  "
$ws.Range("D24").Value = "10.75"
$ws.Range("E24").Value = "  -3.15%  "
$ws.Range("D25").Value = "2.051.40"
$ws.Range("E25").Value = "  -4.00%  "
$ws.Range("D26").Value = "2.038"
$ws.Range("E26").Value = "  -0.85%  "
$ws.Range("D27").Value = "151.41"
$ws.Range("E27").Value = "  -4.31%  "
$ws.Range("D28").Value = "18.18"
$ws.Range("E28").Value = "  -4.51%  "
$ws.Range("D29").Value = "5.319"
$ws.Range("E29").Value = "  -5.39%  "
$ws.Range("D30").Value = "111.10"
$ws.Range("E30").Value = "  -5.52%  "
$ws.Range("D31").Value = "1.727"
$ws.Range("E31").Value = "  -6.25%  "
$ws.Range("D32").Value = "0.08873"
$ws.Range("E32").Value = "  -4.74%  "
$ws.Range("D33").Value = "0.7754"
$ws.Range("E33").Value = "  -10.39%  "
$ws.Range("D34").Value = "4.471"
$ws.Range("E34").Value = "  -12.50%  "
$ws.Range("D35").Value = "2.861"
$ws.Range("E35").Value = "  -5.16%  "
$ws.Range("D36").Value = "1.082"
$ws.Range("E36").Value = "  -13.15%  "
$ws.Range("D37").Value = "0.9997"
$ws.Range("E37").Value = "  -0.42%  "
$ws.Range("D38").Value = "0.05396"
$ws.Range("E38").Value = "  -5.23%  "
$ws.Range("D39").Value = "1.096"
$ws.Range("E39").Value = "  -4.90%  "
$ws.Range("D40").Value = "2.938"
$ws.Range("E40").Value = "  -5.36%  "
$ws.Range("E41").Value = "  -6.81%  "
$ws.Range("D42").Value = "0.5052"
$ws.Range("E42").Value = "  -8.27%  "
$ws.Range("D43").Value = "6.786"
$ws.Range("E43").Value = "  -9.08%  "
$ws.Range("D44").Value = "0.1632"
$ws.Range("E44").Value = "  -7.13%  "
$ws.Range("D45").Value = "0.06601"
$ws.Range("E45").Value = "  -4.93%  "
$ws.Range("D46").Value = "8.209"
$ws.Range("E46").Value = "  -12.13%  "
$ws.Range("D47").Value = "106.03"
$ws.Range("E47").Value = "  -4.41%  "
$ws.Range("D48").Value = "0.4699"
$ws.Range("E48").Value = "  -9.20%  "
$ws.Range("D49").Value = "10.16"
$ws.Range("E49").Value = "  -9.39%  "
$ws.Range("D50").Value = "0.9992"
$ws.Range("E50").Value = "  -0.45%  "
$ws.Range("D51").Value = "1.644"
$ws.Range("E51").Value = "  -7.02%  "

# Restore default (unstyled) formatting so cells match the original
# workbook (no explicit style index on data cells).
$ws.Range("D2:E51").ClearFormats()
